$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-14 across columns A (Date), B (Name), C (Price),
# D (Country), E (Vat Number). A blank string in E means the cell should
# have no value at all (cleared / removed).
$data = @(
  @("2021-04-25", "Owain",              255000, "United Kingdom", ""),
  @("2021-04-23", "Owain",              16500,  "United Kingdom", ""),
  @("2021-04-23", "Owain",              15000,  "United Kingdom", ""),
  @("2021-03-22", "example company",    950,    "United Kingdom", "10191882"),
  @("2021-03-22", "different company",  99500,  "United Kingdom", "10195882"),
  @("2021-03-22", "different company",  16000,  "United Kingdom", "10195882"),
  @("2021-03-22", "random company",     2800,   "United Kingdom", "10194882"),
  @("2021-03-22", "new company",        4500,   "United Kingdom", "10131882"),
  @("2021-03-22", "example company",    15000,  "United Kingdom", "10191882"),
  @("2021-03-22", "different company",  30000,  "United Kingdom", "10195882"),
  @("2021-03-22", "random company",     24500,  "United Kingdom", "10194882"),
  @("2021-03-22", "new company",        3950,   "United Kingdom", "10131882"),
  @("2021-03-22", "new company",        2500,   "United Kingdom", "10131882")
)

# The Date column (A) and Vat Number column (E) hold text that Excel would
# otherwise auto-convert to a date serial / plain number, so those cells
# are forced to plain-text format right before being written, then the
# style is put back to Normal so no stray formatting remains.

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 1).Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]

    if ($r[4] -eq "") {
        $ws.Cells.Item($row, 5).Value = ""
    } else {
        $ws.Cells.Item($row, 5).NumberFormat = "@"
        $ws.Cells.Item($row, 5).Value = $r[4]
        $ws.Cells.Item($row, 5).Style = "Normal"
    }
    $row++
}
